$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three catalogoarchivo.madrid.es permalink URLs (shared strings used by C16:C18)
$ws.Range("C16").Value = "http://catalogoarchivo.madrid.es/ms-opac/permalink/4@oai_villa_baratz_es_villa_471497"
$ws.Range("C17").Value = "http://catalogoarchivo.madrid.es/ms-opac/permalink/4@oai_villa_baratz_es_villa_471664"
$ws.Range("C18").Value = "http://catalogoarchivo.madrid.es/ms-opac/permalink/4@oai_villa_baratz_es_villa_471665"

# Touch formatting on the trailing empty cells in column B so they pick up an explicit style
$ws.Range("B31:B36").HorizontalAlignment = -4131

# Leave the active selection on C18, matching where editing left off
$ws.Range("C18").Select()
